# Experiment order generation script
# Regenerates the per-task "task_order" sheets (file-name lists consumed by
# the experiment runner) with a freshly generated set of stim-file orders,
# and renames each tab to match the task whose data it now holds.
#
# Tab *positions* (and therefore the underlying worksheet parts / rIds) do
# not move - only the tab names and their cell contents change:
#   pos1: GNG_TO -> TOL_TO   (grows   5 ->  7 rows)
#   pos2: NB_TO  -> NB_TO    (stays 10 rows, values refreshed)
#   pos3: RS_TO  -> vSAT_TO  (grows   3 ->  5 rows)
#   pos4: TOL_TO -> RS_TO    (shrinks 7 ->  3 rows)
#   pos5: vSAT_TO-> GNG_TO   (stays  5 rows, values refreshed)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Position 1: GNG_TO -> TOL_TO (go/GNG stims -> MM/ZM stims), 5 -> 7 rows
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "TOL_TO-16515889171565597"

$ws1.Range("A6:A7").EntireRow.Insert()
$ws1.Cells.Item(5,1).Copy()
$ws1.Range("A6:A7").PasteSpecial(-4122)

$ws1.Cells.Item(2,1).Value = 0
$ws1.Cells.Item(2,2).Value = "MM_stims-16515889171250255.csv"
$ws1.Cells.Item(3,1).Value = 1
$ws1.Cells.Item(3,2).Value = "ZM_stims-16515889171143084.csv"
$ws1.Cells.Item(4,1).Value = 2
$ws1.Cells.Item(4,2).Value = "MM_stims-16515889171410348.csv"
$ws1.Cells.Item(5,1).Value = 3
$ws1.Cells.Item(5,2).Value = "ZM_stims-16515889171259944.csv"
$ws1.Cells.Item(6,1).Value = 4
$ws1.Cells.Item(6,2).Value = "MM_stims-16515889171565597.csv"
$ws1.Cells.Item(7,1).Value = 5
$ws1.Cells.Item(7,2).Value = "ZM_stims-1651588917142996.csv"

# ---------------------------------------------------------------------
# Position 2: NB_TO -> NB_TO (ZB-match/OB/TB stims refreshed), stays 10 rows
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "NB_TO-16515889198194294"

$ws2.Cells.Item(2,2).Value = "ZB-match_8-16515889172999399.csv"
$ws2.Cells.Item(3,2).Value = "OB-16515889188608818.csv"
$ws2.Cells.Item(4,2).Value = "TB-16515889198079295.csv"
$ws2.Cells.Item(5,2).Value = "OB-16515889185193129.csv"
$ws2.Cells.Item(6,2).Value = "ZB-match_4-1651588917347931.csv"
$ws2.Cells.Item(7,2).Value = "TB-16515889189233277.csv"
$ws2.Cells.Item(8,2).Value = "ZB-match_1-16515889174699535.csv"
$ws2.Cells.Item(9,2).Value = "TB-1651588919266312.csv"
# row 10 (B10) keeps its original value: TB-16512554887287996.csv

# ---------------------------------------------------------------------
# Position 3: RS_TO -> vSAT_TO ("eyes open/closed" -> vSAT/SAT stims), 3 -> 5 rows
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "vSAT_TO-16515889198817332"

$ws3.Range("A4:A5").EntireRow.Insert()
$ws3.Cells.Item(3,1).Copy()
$ws3.Range("A4:A5").PasteSpecial(-4122)

$ws3.Cells.Item(2,1).Value = 0
$ws3.Cells.Item(2,2).Value = "vSAT_stims-16515889198504512.csv"
$ws3.Cells.Item(3,1).Value = 1
$ws3.Cells.Item(3,2).Value = "vSAT_stims-16515889198661945.csv"
$ws3.Cells.Item(4,1).Value = 2
$ws3.Cells.Item(4,2).Value = "SAT_stims-1651588919823452.csv"
$ws3.Cells.Item(5,1).Value = 3
$ws3.Cells.Item(5,2).Value = "SAT_stims-16515889198360057.csv"

# ---------------------------------------------------------------------
# Position 4: TOL_TO -> RS_TO (MM/ZM stims -> "eyes open/closed"), 7 -> 3 rows
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Name = "RS_TO-16515889198836997"

$ws4.Range("A4:A7").EntireRow.Delete()

$ws4.Cells.Item(2,1).Value = 0
$ws4.Cells.Item(2,2).Value = "eyes open"
$ws4.Cells.Item(3,1).Value = 1
$ws4.Cells.Item(3,2).Value = "eyes closed"

# ---------------------------------------------------------------------
# Position 5: vSAT_TO -> GNG_TO (vSAT/SAT stims -> go/GNG stims), stays 5 rows
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Name = "GNG_TO-16515889199304583"

$ws5.Cells.Item(2,2).Value = "go_stims-1651588919884939.csv"
$ws5.Cells.Item(3,2).Value = "GNG_stims-16515889199122446.csv"
$ws5.Cells.Item(4,2).Value = "go_stims-16515889199133317.csv"
$ws5.Cells.Item(5,2).Value = "GNG_stims-1651588919929421.csv"

Write-Output "experiment order regeneration complete"
